$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: reuse "Stack " category (A17 style), and "Question" style (B17), and plain style (C17)
$ws.Range("A17:C17").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)

$ws.Range("A17:C17").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)

$ws.Range("A18").Value = "Stack "
$ws.Range("B18").Value = "150. Evaluate Reverse Polish Notation"
$ws.Range("C18").Value = "Use stack ,iter over tokens, put all the number inside the stack, when u encounter a operator (*,-,+,/) pop 2 nums as num1, num2 and do the specified operations on it and push the result back into the stack,, continue iterating over tokens"

$ws.Range("A19").Value = "Backtracking "
$ws.Range("B19").Value = "22. Generate Parentheses"
$ws.Range("C19").Value = "Do a recursive backtracking. U can add ""("" if num of openParenthes<n,, U can add "")"" to str if num of closing parantheses < num of open Parentheses,, Base case is when open Parentheses = closing Parenthses = n (num of partheses which we can use to genrate valid permutation, this is given)"

$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(19).AutoFit()

$ws.Range("B19").Select()
